$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from Sheet1 to HomePage
$ws.Name = "HomePage"

# Populate the data grid
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Demo Web Shop"
$ws.Range("B2").Value = "Thank you for signing up! A verification email has been sent. We appreciate your interest."
$ws.Range("A2").Value = "SubscribeMessage"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 19.21875
$ws.Columns.Item(2).ColumnWidth = 73.109375

# Selection / view state
$ws.Range("A10").Select() | Out-Null

# Page setup
$ws.PageSetup.Orientation = 1
